$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete the blank rows 5:8 which shifts the rows below (9-35) up to (5-31)
$ws.Range("5:8").Delete()

# Update the selected cell to match the final state
$ws.Range("H6").Select()
